$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking price string must be forced to Text
# format first, otherwise Excel auto-converts them to a Number type.
$textCells = @('D5','D6','D10','D11','D18','D21','D27','D31','D36','D37','D38','D40','D41','D44','D45','D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume %) scraped by the Action.

# Row 2
$ws.Range("D2").Value = '64.356.15'
$ws.Range("E2").Value = '  -0.01%  '

# Row 3
$ws.Range("D3").Value = '3.510.64'
$ws.Range("E3").Value = '  +0.19%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '591.29'
$ws.Range("E5").Value = '  +0.86%  '

# Row 6
$ws.Range("D6").Value = '134.48'
$ws.Range("E6").Value = '  -0.88%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("E9").Value = '  +5.74%  '

# Row 10
$ws.Range("D10").Value = '0.125'
$ws.Range("E10").Value = '  +0.23%  '

# Row 11
$ws.Range("D11").Value = '0.390'
$ws.Range("E11").Value = '  +3.81%  '

# Row 12
$ws.Range("D12").Value = '4.108.56'
$ws.Range("E12").Value = '  +0.22%  '

# Row 13
$ws.Range("E13").Value = '  +1.27%  '

# Row 14
$ws.Range("E14").Value = '  +0.56%  '

# Row 15
$ws.Range("D15").Value = '3.509.17'
$ws.Range("E15").Value = '  +0.19%  '

# Row 16
$ws.Range("E16").Value = '  +1.35%  '

# Row 17
$ws.Range("D17").Value = '64.334.27'
$ws.Range("E17").Value = '  -0.03%  '

# Row 18
$ws.Range("D18").Value = '9.99'
$ws.Range("E18").Value = '  +2.00%  '

# Row 19
$ws.Range("E19").Value = '  +3.12%  '

# Row 21
$ws.Range("D21").Value = '395.02'
$ws.Range("E21").Value = '  +2.89%  '

# Row 22
$ws.Range("E22").Value = '  +1.06%  '

# Row 23
$ws.Range("D23").Value = '3.650.93'

# Row 24
$ws.Range("E24").Value = '  +0.91%  '

# Row 25
$ws.Range("E25").Value = '  -0.01%  '

# Row 26
$ws.Range("E26").Value = '  +0.26%  '

# Row 27
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").Value = '  +2.52%  '

# Row 28
$ws.Range("E28").Value = '  +0.06%  '

# Row 29
$ws.Range("E29").Value = '  -1.87%  '

# Row 30
$ws.Range("E30").Value = '  +1.25%  '

# Row 31
$ws.Range("D31").Value = '8.29'
$ws.Range("E31").Value = '  +0.00%  '

# Row 32
$ws.Range("E32").Value = '  -6.54%  '

# Row 33
$ws.Range("E33").Value = '  +6.11%  '

# Row 34
$ws.Range("D34").Value = '3.540.03'
$ws.Range("E34").Value = '  +0.44%  '

# Row 35
$ws.Range("E35").Value = '  +0.03%  '

# Row 36
$ws.Range("D36").Value = '23.37'
$ws.Range("E36").Value = '  -0.82%  '

# Row 37
$ws.Range("D37").Value = '5.35'
$ws.Range("E37").Value = '  +0.57%  '

# Row 38
$ws.Range("D38").Value = '6.97'
$ws.Range("E38").Value = '  +1.57%  '

# Row 39
$ws.Range("E39").Value = '  +0.22%  '

# Row 40
$ws.Range("D40").Value = '167.21'
$ws.Range("E40").Value = '  +2.10%  '

# Row 41
$ws.Range("D41").Value = '0.0789'

# Row 42
$ws.Range("E42").Value = '  +0.58%  '

# Row 43
$ws.Range("E43").Value = '  -0.03%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '25.20'
$ws.Range("E44").Value = '  -2.51%  '

# Row 45
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '4.44'
$ws.Range("E45").Value = '  +0.82%  '

# Row 46
$ws.Range("E46").Value = '  +1.03%  '

# Row 47
$ws.Range("E47").Value = '  -3.22%  '

# Row 48
$ws.Range("D48").Value = '6.82'
$ws.Range("E48").Value = '  +0.57%  '

# Row 49
$ws.Range("D49").Value = '2.383.69'
$ws.Range("E49").Value = '  -3.62%  '

# Row 50
$ws.Range("E50").Value = '  -1.93%  '

# Row 51
$ws.Range("E51").Value = '  -0.01%  '
